$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Data")

$ws.Range("E2").Value = 0.00191538
$ws.Range("F2").Value = 0.02893401
$ws.Range("G2").Value = 0.0023574651136363638

$ws.Range("E3").Value = 0.00461511
$ws.Range("F3").Value = 0.01011735
$ws.Range("G3").Value = 0.005906208510638298
